$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Jag1"
$ws.Cells.Item(2, 3).Value = "Notch3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 18.78268133333334
$ws.Cells.Item(2, 8).Value = 56.34804400000001
$ws.Cells.Item(2, 9).Value = 0.286679008418643
$ws.Cells.Item(2, 10).Value = 0.286679008418643
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 3.884573333333333
$ws.Cells.Item(2, 14).Value = 11.65372
$ws.Cells.Item(2, 15).Value = 0.09450909866970746
$ws.Cells.Item(2, 16).Value = 0.09450909866970746
$ws.Cells.Item(2, 17).Value = 72.96270303596445
$ws.Cells.Item(2, 18).Value = 656.6643273236801
$ws.Cells.Item(2, 19).Value = 0.02709377469317142
$ws.Cells.Item(2, 20).Value = 0.02709377469317142

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Jag1"
$ws.Cells.Item(3, 3).Value = "Notch3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 18.78268133333334
$ws.Cells.Item(3, 8).Value = 56.34804400000001
$ws.Cells.Item(3, 9).Value = 0.286679008418643
$ws.Cells.Item(3, 10).Value = 0.286679008418643
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 6.662909
$ws.Cells.Item(3, 14).Value = 19.988727
$ws.Cells.Item(3, 15).Value = 0.1621041669376684
$ws.Cells.Item(3, 16).Value = 0.1621041669376685
$ws.Cells.Item(3, 17).Value = 125.1472964999987
$ws.Cells.Item(3, 18).Value = 1126.325668499988
$ws.Cells.Item(3, 19).Value = 0.04647186183822095
$ws.Cells.Item(3, 20).Value = 0.04647186183822096

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Jag1"
$ws.Cells.Item(4, 3).Value = "Notch3"
$ws.Cells.Item(4, 4).Value = "M1"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 18.78268133333334
$ws.Cells.Item(4, 8).Value = 56.34804400000001
$ws.Cells.Item(4, 9).Value = 0.286679008418643
$ws.Cells.Item(4, 10).Value = 0.286679008418643
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.04919333333333333
$ws.Cells.Item(4, 14).Value = 0.14758
$ws.Cells.Item(4, 15).Value = 0.001196841247402154
$ws.Cells.Item(4, 16).Value = 0.001196841247402154
$ws.Cells.Item(4, 17).Value = 0.9239827037244446
$ws.Cells.Item(4, 18).Value = 8.315844333520001
$ws.Cells.Item(4, 19).Value = 0.0003431092620397811
$ws.Cells.Item(4, 20).Value = 0.0003431092620397811

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Jag1"
$ws.Cells.Item(5, 3).Value = "Notch3"
$ws.Cells.Item(5, 4).Value = "M2"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 18.78268133333334
$ws.Cells.Item(5, 8).Value = 56.34804400000001
$ws.Cells.Item(5, 9).Value = 0.286679008418643
$ws.Cells.Item(5, 10).Value = 0.286679008418643
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.106491
$ws.Cells.Item(5, 14).Value = 0.319473
$ws.Cells.Item(5, 15).Value = 0.002590855561941376
$ws.Cells.Item(5, 16).Value = 0.002590855561941376
$ws.Cells.Item(5, 17).Value = 2.000186517868
$ws.Cells.Item(5, 18).Value = 18.001678660812
$ws.Cells.Item(5, 19).Value = 0.0007427439034532797
$ws.Cells.Item(5, 20).Value = 0.0007427439034532798

# Row 6
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Jag1"
$ws.Cells.Item(6, 3).Value = "Notch3"
$ws.Cells.Item(6, 4).Value = "sCs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 18.78268133333334
$ws.Cells.Item(6, 8).Value = 56.34804400000001
$ws.Cells.Item(6, 9).Value = 0.286679008418643
$ws.Cells.Item(6, 10).Value = 0.286679008418643
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 30.399472
$ws.Cells.Item(6, 14).Value = 91.19841600000001
$ws.Cells.Item(6, 15).Value = 0.7395990375832805
$ws.Cells.Item(6, 16).Value = 0.7395990375832806
$ws.Cells.Item(6, 17).Value = 570.9835952775895
$ws.Cells.Item(6, 18).Value = 5138.852357498306
$ws.Cells.Item(6, 19).Value = 0.2120275187217575
$ws.Cells.Item(6, 20).Value = 0.2120275187217575

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Jag1"
$ws.Cells.Item(7, 3).Value = "Notch3"
$ws.Cells.Item(7, 4).Value = "ECs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 14.452944
$ws.Cells.Item(7, 8).Value = 43.358832
$ws.Cells.Item(7, 9).Value = 0.2205944711044544
$ws.Cells.Item(7, 10).Value = 0.2205944711044544
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 3.884573333333333
$ws.Cells.Item(7, 14).Value = 11.65372
$ws.Cells.Item(7, 15).Value = 0.09450909866970746
$ws.Cells.Item(7, 16).Value = 0.09450909866970746
$ws.Cells.Item(7, 17).Value = 56.14352085056
$ws.Cells.Item(7, 18).Value = 505.29168765504
$ws.Cells.Item(7, 19).Value = 0.02084818463560281
$ws.Cells.Item(7, 20).Value = 0.02084818463560281

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Jag1"
$ws.Cells.Item(8, 3).Value = "Notch3"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 14.452944
$ws.Cells.Item(8, 8).Value = 43.358832
$ws.Cells.Item(8, 9).Value = 0.2205944711044544
$ws.Cells.Item(8, 10).Value = 0.2205944711044544
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 6.662909
$ws.Cells.Item(8, 14).Value = 19.988727
$ws.Cells.Item(8, 15).Value = 0.1621041669376684
$ws.Cells.Item(8, 16).Value = 0.1621041669376685
$ws.Cells.Item(8, 17).Value = 96.29865065409601
$ws.Cells.Item(8, 18).Value = 866.687855886864
$ws.Cells.Item(8, 19).Value = 0.03575928296944315
$ws.Cells.Item(8, 20).Value = 0.03575928296944315

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Jag1"
$ws.Cells.Item(9, 3).Value = "Notch3"
$ws.Cells.Item(9, 4).Value = "M1"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 14.452944
$ws.Cells.Item(9, 8).Value = 43.358832
$ws.Cells.Item(9, 9).Value = 0.2205944711044544
$ws.Cells.Item(9, 10).Value = 0.2205944711044544
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.04919333333333333
$ws.Cells.Item(9, 14).Value = 0.14758
$ws.Cells.Item(9, 15).Value = 0.001196841247402154
$ws.Cells.Item(9, 16).Value = 0.001196841247402154
$ws.Cells.Item(9, 17).Value = 0.71098849184
$ws.Cells.Item(9, 18).Value = 6.398896426559999
$ws.Cells.Item(9, 19).Value = 0.0002640165619666735
$ws.Cells.Item(9, 20).Value = 0.0002640165619666735

# Row 10
$ws.Cells.Item(10, 1).Value = "FAPs"
$ws.Cells.Item(10, 2).Value = "Jag1"
$ws.Cells.Item(10, 3).Value = "Notch3"
$ws.Cells.Item(10, 4).Value = "M2"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 14.452944
$ws.Cells.Item(10, 8).Value = 43.358832
$ws.Cells.Item(10, 9).Value = 0.2205944711044544
$ws.Cells.Item(10, 10).Value = 0.2205944711044544
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.106491
$ws.Cells.Item(10, 14).Value = 0.319473
$ws.Cells.Item(10, 15).Value = 0.002590855561941376
$ws.Cells.Item(10, 16).Value = 0.002590855561941376
$ws.Cells.Item(10, 17).Value = 1.539108459504
$ws.Cells.Item(10, 18).Value = 13.851976135536
$ws.Cells.Item(10, 19).Value = 0.0005715284123944917
$ws.Cells.Item(10, 20).Value = 0.0005715284123944918

# Row 11
$ws.Cells.Item(11, 1).Value = "FAPs"
$ws.Cells.Item(11, 2).Value = "Jag1"
$ws.Cells.Item(11, 3).Value = "Notch3"
$ws.Cells.Item(11, 4).Value = "sCs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 14.452944
$ws.Cells.Item(11, 8).Value = 43.358832
$ws.Cells.Item(11, 9).Value = 0.2205944711044544
$ws.Cells.Item(11, 10).Value = 0.2205944711044544
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 30.399472
$ws.Cells.Item(11, 14).Value = 91.19841600000001
$ws.Cells.Item(11, 15).Value = 0.7395990375832805
$ws.Cells.Item(11, 16).Value = 0.7395990375832806
$ws.Cells.Item(11, 17).Value = 439.3618664455681
$ws.Cells.Item(11, 18).Value = 3954.256798010113
$ws.Cells.Item(11, 19).Value = 0.1631514585250472
$ws.Cells.Item(11, 20).Value = 0.1631514585250473

# Row 12
$ws.Cells.Item(12, 1).Value = "M1"
$ws.Cells.Item(12, 2).Value = "Jag1"
$ws.Cells.Item(12, 3).Value = "Notch3"
$ws.Cells.Item(12, 4).Value = "ECs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 20.02445833333333
$ws.Cells.Item(12, 8).Value = 60.073375
$ws.Cells.Item(12, 9).Value = 0.3056321809034097
$ws.Cells.Item(12, 10).Value = 0.3056321809034098
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 3.884573333333333
$ws.Cells.Item(12, 14).Value = 11.65372
$ws.Cells.Item(12, 15).Value = 0.09450909866970746
$ws.Cells.Item(12, 16).Value = 0.09450909866970746
$ws.Cells.Item(12, 17).Value = 77.78647685611111
$ws.Cells.Item(12, 18).Value = 700.078291705
$ws.Cells.Item(12, 19).Value = 0.02888502194163823
$ws.Cells.Item(12, 20).Value = 0.02888502194163823

# Row 13
$ws.Cells.Item(13, 1).Value = "M1"
$ws.Cells.Item(13, 2).Value = "Jag1"
$ws.Cells.Item(13, 3).Value = "Notch3"
$ws.Cells.Item(13, 4).Value = "FAPs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 20.02445833333333
$ws.Cells.Item(13, 8).Value = 60.073375
$ws.Cells.Item(13, 9).Value = 0.3056321809034097
$ws.Cells.Item(13, 10).Value = 0.3056321809034098
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 6.662909
$ws.Cells.Item(13, 14).Value = 19.988727
$ws.Cells.Item(13, 15).Value = 0.1621041669376684
$ws.Cells.Item(13, 16).Value = 0.1621041669376685
$ws.Cells.Item(13, 17).Value = 133.4211436492917
$ws.Cells.Item(13, 18).Value = 1200.790292843625
$ws.Cells.Item(13, 19).Value = 0.04954425007469001
$ws.Cells.Item(13, 20).Value = 0.04954425007469002

# Row 14
$ws.Cells.Item(14, 1).Value = "M1"
$ws.Cells.Item(14, 2).Value = "Jag1"
$ws.Cells.Item(14, 3).Value = "Notch3"
$ws.Cells.Item(14, 4).Value = "M1"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 20.02445833333333
$ws.Cells.Item(14, 8).Value = 60.073375
$ws.Cells.Item(14, 9).Value = 0.3056321809034097
$ws.Cells.Item(14, 10).Value = 0.3056321809034098
$ws.Cells.Item(14, 11).Value = 1
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 0.04919333333333333
$ws.Cells.Item(14, 14).Value = 0.14758
$ws.Cells.Item(14, 15).Value = 0.001196841247402154
$ws.Cells.Item(14, 16).Value = 0.001196841247402154
$ws.Cells.Item(14, 17).Value = 0.985069853611111
$ws.Cells.Item(14, 18).Value = 8.865628682499999
$ws.Cells.Item(14, 19).Value = 0.0003657932006386775
$ws.Cells.Item(14, 20).Value = 0.0003657932006386776

# Row 15
$ws.Cells.Item(15, 1).Value = "M1"
$ws.Cells.Item(15, 2).Value = "Jag1"
$ws.Cells.Item(15, 3).Value = "Notch3"
$ws.Cells.Item(15, 4).Value = "M2"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 20.02445833333333
$ws.Cells.Item(15, 8).Value = 60.073375
$ws.Cells.Item(15, 9).Value = 0.3056321809034097
$ws.Cells.Item(15, 10).Value = 0.3056321809034098
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 0.106491
$ws.Cells.Item(15, 14).Value = 0.319473
$ws.Cells.Item(15, 15).Value = 0.002590855561941376
$ws.Cells.Item(15, 16).Value = 0.002590855561941376
$ws.Cells.Item(15, 17).Value = 2.132424592375
$ws.Cells.Item(15, 18).Value = 19.191821331375
$ws.Cells.Item(15, 19).Value = 0.0007918488358018718
$ws.Cells.Item(15, 20).Value = 0.0007918488358018721

# Row 16
$ws.Cells.Item(16, 1).Value = "M1"
$ws.Cells.Item(16, 2).Value = "Jag1"
$ws.Cells.Item(16, 3).Value = "Notch3"
$ws.Cells.Item(16, 4).Value = "sCs"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 20.02445833333333
$ws.Cells.Item(16, 8).Value = 60.073375
$ws.Cells.Item(16, 9).Value = 0.3056321809034097
$ws.Cells.Item(16, 10).Value = 0.3056321809034098
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 30.399472
$ws.Cells.Item(16, 14).Value = 91.19841600000001
$ws.Cells.Item(16, 15).Value = 0.7395990375832805
$ws.Cells.Item(16, 16).Value = 0.7395990375832806
$ws.Cells.Item(16, 17).Value = 608.7329604193334
$ws.Cells.Item(16, 18).Value = 5478.596643774001
$ws.Cells.Item(16, 19).Value = 0.2260452668506409
$ws.Cells.Item(16, 20).Value = 0.226045266850641

# Row 17
$ws.Cells.Item(17, 1).Value = "M2"
$ws.Cells.Item(17, 2).Value = "Jag1"
$ws.Cells.Item(17, 3).Value = "Notch3"
$ws.Cells.Item(17, 4).Value = "ECs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 6.777317666666666
$ws.Cells.Item(17, 8).Value = 20.331953
$ws.Cells.Item(17, 9).Value = 0.1034418182333792
$ws.Cells.Item(17, 10).Value = 0.1034418182333792
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 3.884573333333333
$ws.Cells.Item(17, 14).Value = 11.65372
$ws.Cells.Item(17, 15).Value = 0.09450909866970746
$ws.Cells.Item(17, 16).Value = 0.09450909866970746
$ws.Cells.Item(17, 17).Value = 26.32698747946222
$ws.Cells.Item(17, 18).Value = 236.94288731516
$ws.Cells.Item(17, 19).Value = 0.009776193005992374
$ws.Cells.Item(17, 20).Value = 0.009776193005992376

# Row 18
$ws.Cells.Item(18, 1).Value = "M2"
$ws.Cells.Item(18, 2).Value = "Jag1"
$ws.Cells.Item(18, 3).Value = "Notch3"
$ws.Cells.Item(18, 4).Value = "FAPs"
$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 6.777317666666666
$ws.Cells.Item(18, 8).Value = 20.331953
$ws.Cells.Item(18, 9).Value = 0.1034418182333792
$ws.Cells.Item(18, 10).Value = 0.1034418182333792
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 6.662909
$ws.Cells.Item(18, 14).Value = 19.988727
$ws.Cells.Item(18, 15).Value = 0.1621041669376684
$ws.Cells.Item(18, 16).Value = 0.1621041669376685
$ws.Cells.Item(18, 17).Value = 45.15665087709232
$ws.Cells.Item(18, 18).Value = 406.409857893831
$ws.Cells.Item(18, 19).Value = 0.01676834977123965
$ws.Cells.Item(18, 20).Value = 0.01676834977123965

# Row 19
$ws.Cells.Item(19, 1).Value = "M2"
$ws.Cells.Item(19, 2).Value = "Jag1"
$ws.Cells.Item(19, 3).Value = "Notch3"
$ws.Cells.Item(19, 4).Value = "M1"
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 6.777317666666666
$ws.Cells.Item(19, 8).Value = 20.331953
$ws.Cells.Item(19, 9).Value = 0.1034418182333792
$ws.Cells.Item(19, 10).Value = 0.1034418182333792
$ws.Cells.Item(19, 11).Value = 1
$ws.Cells.Item(19, 12).Value = 0.3333333333333333
$ws.Cells.Item(19, 13).Value = 0.04919333333333333
$ws.Cells.Item(19, 14).Value = 0.14758
$ws.Cells.Item(19, 15).Value = 0.001196841247402154
$ws.Cells.Item(19, 16).Value = 0.001196841247402154
$ws.Cells.Item(19, 17).Value = 0.3333988470822222
$ws.Cells.Item(19, 18).Value = 3.000589623739999
$ws.Cells.Item(19, 19).Value = 0.0001238034347679843
$ws.Cells.Item(19, 20).Value = 0.0001238034347679844

# Row 20
$ws.Cells.Item(20, 1).Value = "M2"
$ws.Cells.Item(20, 2).Value = "Jag1"
$ws.Cells.Item(20, 3).Value = "Notch3"
$ws.Cells.Item(20, 4).Value = "M2"
$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 6.777317666666666
$ws.Cells.Item(20, 8).Value = 20.331953
$ws.Cells.Item(20, 9).Value = 0.1034418182333792
$ws.Cells.Item(20, 10).Value = 0.1034418182333792
$ws.Cells.Item(20, 11).Value = 3
$ws.Cells.Item(20, 12).Value = 1
$ws.Cells.Item(20, 13).Value = 0.106491
$ws.Cells.Item(20, 14).Value = 0.319473
$ws.Cells.Item(20, 15).Value = 0.002590855561941376
$ws.Cells.Item(20, 16).Value = 0.002590855561941376
$ws.Cells.Item(20, 17).Value = 0.7217233356409999
$ws.Cells.Item(20, 18).Value = 6.495510020768999
$ws.Cells.Item(20, 19).Value = 0.0002680028101072792
$ws.Cells.Item(20, 20).Value = 0.0002680028101072793

# Row 21
$ws.Cells.Item(21, 1).Value = "M2"
$ws.Cells.Item(21, 2).Value = "Jag1"
$ws.Cells.Item(21, 3).Value = "Notch3"
$ws.Cells.Item(21, 4).Value = "sCs"
$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = 6.777317666666666
$ws.Cells.Item(21, 8).Value = 20.331953
$ws.Cells.Item(21, 9).Value = 0.1034418182333792
$ws.Cells.Item(21, 10).Value = 0.1034418182333792
$ws.Cells.Item(21, 11).Value = 3
$ws.Cells.Item(21, 12).Value = 1
$ws.Cells.Item(21, 13).Value = 30.399472
$ws.Cells.Item(21, 14).Value = 91.19841600000001
$ws.Cells.Item(21, 15).Value = 0.7395990375832805
$ws.Cells.Item(21, 16).Value = 0.7395990375832806
$ws.Cells.Item(21, 17).Value = 206.0268786429387
$ws.Cells.Item(21, 18).Value = 1854.241907786448
$ws.Cells.Item(21, 19).Value = 0.07650546921127185
$ws.Cells.Item(21, 20).Value = 0.07650546921127188

# Row 22
$ws.Cells.Item(22, 1).Value = "sCs"
$ws.Cells.Item(22, 2).Value = "Jag1"
$ws.Cells.Item(22, 3).Value = "Notch3"
$ws.Cells.Item(22, 4).Value = "ECs"
$ws.Cells.Item(22, 5).Value = 3
$ws.Cells.Item(22, 6).Value = 1
$ws.Cells.Item(22, 7).Value = 5.480759333333334
$ws.Cells.Item(22, 8).Value = 16.442278
$ws.Cells.Item(22, 9).Value = 0.08365252134011374
$ws.Cells.Item(22, 10).Value = 0.08365252134011374
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 12).Value = 1
$ws.Cells.Item(22, 13).Value = 3.884573333333333
$ws.Cells.Item(22, 14).Value = 11.65372
$ws.Cells.Item(22, 15).Value = 0.09450909866970746
$ws.Cells.Item(22, 16).Value = 0.09450909866970746
$ws.Cells.Item(22, 17).Value = 21.29041155268445
$ws.Cells.Item(22, 18).Value = 191.61370397416
$ws.Cells.Item(22, 19).Value = 0.007905924393302618
$ws.Cells.Item(22, 20).Value = 0.007905924393302618

# Row 23
$ws.Cells.Item(23, 1).Value = "sCs"
$ws.Cells.Item(23, 2).Value = "Jag1"
$ws.Cells.Item(23, 3).Value = "Notch3"
$ws.Cells.Item(23, 4).Value = "FAPs"
$ws.Cells.Item(23, 5).Value = 3
$ws.Cells.Item(23, 6).Value = 1
$ws.Cells.Item(23, 7).Value = 5.480759333333334
$ws.Cells.Item(23, 8).Value = 16.442278
$ws.Cells.Item(23, 9).Value = 0.08365252134011374
$ws.Cells.Item(23, 10).Value = 0.08365252134011374
$ws.Cells.Item(23, 11).Value = 3
$ws.Cells.Item(23, 12).Value = 1
$ws.Cells.Item(23, 13).Value = 6.662909
$ws.Cells.Item(23, 14).Value = 19.988727
$ws.Cells.Item(23, 15).Value = 0.1621041669376684
$ws.Cells.Item(23, 16).Value = 0.1621041669376685
$ws.Cells.Item(23, 17).Value = 36.51780068890067
$ws.Cells.Item(23, 18).Value = 328.6602062001061
$ws.Cells.Item(23, 19).Value = 0.01356042228407467
$ws.Cells.Item(23, 20).Value = 0.01356042228407467

# Row 24
$ws.Cells.Item(24, 1).Value = "sCs"
$ws.Cells.Item(24, 2).Value = "Jag1"
$ws.Cells.Item(24, 3).Value = "Notch3"
$ws.Cells.Item(24, 4).Value = "M1"
$ws.Cells.Item(24, 5).Value = 3
$ws.Cells.Item(24, 6).Value = 1
$ws.Cells.Item(24, 7).Value = 5.480759333333334
$ws.Cells.Item(24, 8).Value = 16.442278
$ws.Cells.Item(24, 9).Value = 0.08365252134011374
$ws.Cells.Item(24, 10).Value = 0.08365252134011374
$ws.Cells.Item(24, 11).Value = 1
$ws.Cells.Item(24, 12).Value = 0.3333333333333333
$ws.Cells.Item(24, 13).Value = 0.04919333333333333
$ws.Cells.Item(24, 14).Value = 0.14758
$ws.Cells.Item(24, 15).Value = 0.001196841247402154
$ws.Cells.Item(24, 16).Value = 0.001196841247402154
$ws.Cells.Item(24, 17).Value = 0.2696168208044445
$ws.Cells.Item(24, 18).Value = 2.42655138724
$ws.Cells.Item(24, 19).Value = 0.000100118787989037
$ws.Cells.Item(24, 20).Value = 0.000100118787989037

# Row 25
$ws.Cells.Item(25, 1).Value = "sCs"
$ws.Cells.Item(25, 2).Value = "Jag1"
$ws.Cells.Item(25, 3).Value = "Notch3"
$ws.Cells.Item(25, 4).Value = "M2"
$ws.Cells.Item(25, 5).Value = 3
$ws.Cells.Item(25, 6).Value = 1
$ws.Cells.Item(25, 7).Value = 5.480759333333334
$ws.Cells.Item(25, 8).Value = 16.442278
$ws.Cells.Item(25, 9).Value = 0.08365252134011374
$ws.Cells.Item(25, 10).Value = 0.08365252134011374
$ws.Cells.Item(25, 11).Value = 3
$ws.Cells.Item(25, 12).Value = 1
$ws.Cells.Item(25, 13).Value = 0.106491
$ws.Cells.Item(25, 14).Value = 0.319473
$ws.Cells.Item(25, 15).Value = 0.002590855561941376
$ws.Cells.Item(25, 16).Value = 0.002590855561941376
$ws.Cells.Item(25, 17).Value = 0.5836515421660001
$ws.Cells.Item(25, 18).Value = 5.252863879494001
$ws.Cells.Item(25, 19).Value = 0.0002167316001844533
$ws.Cells.Item(25, 20).Value = 0.0002167316001844533

# Row 26
$ws.Cells.Item(26, 1).Value = "sCs"
$ws.Cells.Item(26, 2).Value = "Jag1"
$ws.Cells.Item(26, 3).Value = "Notch3"
$ws.Cells.Item(26, 4).Value = "sCs"
$ws.Cells.Item(26, 5).Value = 3
$ws.Cells.Item(26, 6).Value = 1
$ws.Cells.Item(26, 7).Value = 5.480759333333334
$ws.Cells.Item(26, 8).Value = 16.442278
$ws.Cells.Item(26, 9).Value = 0.08365252134011374
$ws.Cells.Item(26, 10).Value = 0.08365252134011374
$ws.Cells.Item(26, 11).Value = 3
$ws.Cells.Item(26, 12).Value = 1
$ws.Cells.Item(26, 13).Value = 30.399472
$ws.Cells.Item(26, 14).Value = 91.19841600000001
$ws.Cells.Item(26, 15).Value = 0.7395990375832805
$ws.Cells.Item(26, 16).Value = 0.7395990375832806
$ws.Cells.Item(26, 17).Value = 166.6121898924054
$ws.Cells.Item(26, 18).Value = 1499.509709031648
$ws.Cells.Item(26, 19).Value = 0.06186932427456295
$ws.Cells.Item(26, 20).Value = 0.06186932427456296
